# Update the Status column on Sheet1 for the rows whose sign-up test
# results flipped from "Fail" to "Pass".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(4, 12, 13, 14, 16)
foreach ($r in $rows) {
    $ws.Range("C$r").Value = "Pass"
}
